# hotfix event
# Update the "achievement" worksheet:
#   - extend the AND-chained condition formula for achievement #219 (row 121)
#     with one more ATLT clause
#   - append three new achievements (rows 165-167), reusing the two existing
#     blank spacer rows (165, 166) and growing one brand-new data row (167)
#   - re-create a trailing blank spacer row (168)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 121 (achievement id 219): extend the concatenated condition ------
$ws.Cells.Item(121, 5).Value = "(ATLT?[1023])&(ATLT?[1048])&(ATLT?[1064])&(ATLT?[1114])&(ATLT?[1135])&(ATLT?[1141])&(ATLT?[1147])"

# --- New row 165: achievement id 263 (reuses existing blank row's format) -
$ws.Cells.Item(165, 1).Value = 263
$ws.Cells.Item(165, 2).Value = "经 典 老 歌"
$ws.Cells.Item(165, 3).Value = "自 投 罗 网"
$ws.Cells.Item(165, 4).Value = 2
$ws.Cells.Item(165, 5).Value = "TLT?[2036]"
$ws.Cells.Item(165, 6).Value = 1
$ws.Cells.Item(165, 7).Value = "START"

# --- New row 166: achievement id 264 (reuses existing blank row's format) -
$ws.Cells.Item(166, 1).Value = 264
$ws.Cells.Item(166, 2).Value = "莎比"
$ws.Cells.Item(166, 3).Value = "集齐四大悲剧"
$ws.Cells.Item(166, 4).Value = 2
$ws.Cells.Item(166, 5).Value = "(ATLT?[2028])&(ATLT?[2029])&(ATLT?[2030])&(ATLT?[2031])"
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = "START"

# --- New row 167: achievement id 265 (brand-new row, copy formatting first)
$ws.Range("A164:G164").Copy()
$ws.Range("A167:G167").PasteSpecial(-4122)
$ws.Rows.Item(167).RowHeight = 39.6

$ws.Cells.Item(167, 1).Value = 265
$ws.Cells.Item(167, 2).Value = "死了但没完全死"
$ws.Cells.Item(167, 3).Value = "死而复生"
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = "EVT?[20000,20001,11504]"
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = "TRAJECTORY"

# --- New trailing blank row 168 (mirrors the previous blank spacer rows) --
$ws.Range("E165").Copy()
$ws.Range("E168").PasteSpecial(-4122)
$ws.Rows.Item(168).RowHeight = 39.6

# --- Restore the selection the author left the sheet on -------------------
$ws.Range("G122").Select()
